# Insert 3 new weekly rows of Acelga (Chard) price data at the top of the
# "logica_diaria" history block (row 445), pushing the existing rows
# (445-563) down to (448-566). This mirrors a new week of data being
# prepended ahead of the most recent previously-recorded week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 445 (shifts 445:563 -> 448:566)
$ws.Range("A445:A447").EntireRow.Insert()

# Row 445: Extra
$ws.Cells.Item(445, 1).Value = 6
$ws.Cells.Item(445, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(445, 3).Value = "Metropolitana"
$ws.Cells.Item(445, 4).Value = 44508
$ws.Cells.Item(445, 5).Value = 13
$ws.Cells.Item(445, 6).Value = 100112009
$ws.Cells.Item(445, 7).Value = "Acelga"
$ws.Cells.Item(445, 8).Value = "Sin especificar"
$ws.Cells.Item(445, 9).Value = "Extra"
$ws.Cells.Item(445, 10).Value = 26000
$ws.Cells.Item(445, 11).Value = 11000
$ws.Cells.Item(445, 12).Value = 11000
$ws.Cells.Item(445, 13).Value = 11000
$ws.Cells.Item(445, 14).Value = "`$/docena de atados"
$ws.Cells.Item(445, 15).Value = "Región Metropolitana"
$ws.Cells.Item(445, 16).Value = 3667
$ws.Cells.Item(445, 17).Value = 3
$ws.Cells.Item(445, 18).Value = "Hortaliza"

# Row 446: Primera
$ws.Cells.Item(446, 1).Value = 6
$ws.Cells.Item(446, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(446, 3).Value = "Metropolitana"
$ws.Cells.Item(446, 4).Value = 44508
$ws.Cells.Item(446, 5).Value = 13
$ws.Cells.Item(446, 6).Value = 100112009
$ws.Cells.Item(446, 7).Value = "Acelga"
$ws.Cells.Item(446, 8).Value = "Sin especificar"
$ws.Cells.Item(446, 9).Value = "Primera"
$ws.Cells.Item(446, 10).Value = 32000
$ws.Cells.Item(446, 11).Value = 15000
$ws.Cells.Item(446, 12).Value = 15000
$ws.Cells.Item(446, 13).Value = 15000
$ws.Cells.Item(446, 14).Value = "`$/docena de atados"
$ws.Cells.Item(446, 15).Value = "Región Metropolitana"
$ws.Cells.Item(446, 16).Value = 5000
$ws.Cells.Item(446, 17).Value = 3
$ws.Cells.Item(446, 18).Value = "Hortaliza"

# Row 447: Segunda
$ws.Cells.Item(447, 1).Value = 6
$ws.Cells.Item(447, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(447, 3).Value = "Metropolitana"
$ws.Cells.Item(447, 4).Value = 44508
$ws.Cells.Item(447, 5).Value = 13
$ws.Cells.Item(447, 6).Value = 100112009
$ws.Cells.Item(447, 7).Value = "Acelga"
$ws.Cells.Item(447, 8).Value = "Sin especificar"
$ws.Cells.Item(447, 9).Value = "Segunda"
$ws.Cells.Item(447, 10).Value = 9000
$ws.Cells.Item(447, 11).Value = 12000
$ws.Cells.Item(447, 12).Value = 12000
$ws.Cells.Item(447, 13).Value = 12000
$ws.Cells.Item(447, 14).Value = "`$/docena de atados"
$ws.Cells.Item(447, 15).Value = "Región Metropolitana"
$ws.Cells.Item(447, 16).Value = 4000
$ws.Cells.Item(447, 17).Value = 3
$ws.Cells.Item(447, 18).Value = "Hortaliza"
